# Apply OOXML diff to asyntask_app.docx
$d = $word.ActiveDocument

# --- Change 1: merge "A" + "synctask" + " có 4 " runs into a single run ---
$d.Content.Find.Execute(
    "Asynctask có 4 ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Asynctask có 4 ", 2
) | Out-Null

# --- Change 2: merge "3 tk còn" + " có thể gọi hoặc ko gọi cũng ok" into one run ---
$d.Content.Find.Execute(
    "3 tk còn có thể gọi hoặc ko gọi cũng ok", $true, $false, $false, $false, $false,
    $true, 1, $false, "3 tk còn có thể gọi hoặc ko gọi cũng ok", 2
) | Out-Null

# --- Change 3: merge "doInBackground (...) sau khi " + "doInBackground" + " hoạt động " into one run ---
$d.Content.Find.Execute(
    "doInBackground (hàm sử lý chính công việc) sau khi doInBackground hoạt động ", $true, $false, $false, $false, $false,
    $true, 1, $false, "doInBackground (hàm sử lý chính công việc) sau khi doInBackground hoạt động ", 2
) | Out-Null

# --- Change 4: append the new "Params ..." content into the empty paragraph that
# immediately follows another empty (firstLine-indented) paragraph, right before the
# trailing blank paragraphs at the end of the document. ---
$paras = $d.Paragraphs
$target = $null
for ($i = 2; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    $candLen = $cand.Range.End - $cand.Range.Start
    $prev = $paras.Item($i - 1)
    $prevLen = $prev.Range.End - $prev.Range.Start
    if ($candLen -eq 1 -and $cand.Format.FirstLineIndent -eq 36 -and $cand.Format.LeftIndent -eq 0 `
        -and $prevLen -eq 1 -and $prev.Format.FirstLineIndent -eq 36) {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target empty paragraph for Change 4"
}

$pStart = $target.Range.Start
$marker = "PLACEHOLDER_INSERT_MARK"
$target.Range.InsertAfter($marker)
$ins = $d.Range($pStart, $pStart + $marker.Length)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Params là tham số có thể được nhận từ </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>execute(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">), </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:b/><w:bCs/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>có thể là một mảng các tham số con</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">. và Params sẽ là input của </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>doInBackground(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">), Process là input của onProgressUpdate(), đầu ra này nhận từ doInBackground() thông qua phương thức publishProgress(). </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>onProgressUpdate(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">) có thể cập nhật giao diện lúc runtime. Result là đầu ra của </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>doInBackground(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI"/><w:color w:val="292B2C"/><w:sz w:val="27"/><w:szCs w:val="27"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>) và chính là kết quả trả về ở onPostExecute(). AsyncTask chạy trên Worker Thread còn Handler chạy trên Main Thread (hay Thread sinh ra nó).</w:t></w:r></w:p>'
$ins.InsertXML($xml)
